# Update "苏州-漫展信息.xlsx" per the commit "Update gh-pages to output generated at 456a3b4"
# Two sheets are touched:
#   - 展览   (exhibitions)  : a handful of "want-to-go" (F column) counters bumped up
#   - 全部类型 (all types)   : the same counters bumped, plus one more row added for the
#                             "Come in joy" event (duplicated/re-inserted at row 31 with
#                             its updated counter), pushing every later row down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览": bump the "想去人数" (F column) counters for the rows that changed.
# ---------------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")

$wsExpo.Range("F2").Value  = 614
$wsExpo.Range("F4").Value  = 1288
$wsExpo.Range("F5").Value  = 1153
$wsExpo.Range("F6").Value  = 14289
$wsExpo.Range("F7").Value  = 16323
$wsExpo.Range("F9").Value  = 90
$wsExpo.Range("F18").Value = 99
$wsExpo.Range("F20").Value = 1247
$wsExpo.Range("F23").Value = 33
$wsExpo.Range("F24").Value = 6572
$wsExpo.Range("F25").Value = 969
$wsExpo.Range("F26").Value = 15
$wsExpo.Range("F29").Value = 5707
$wsExpo.Range("F33").Value = 4757

# ---------------------------------------------------------------------------
# Sheet "全部类型": same counter bumps (rows line up slightly differently here
# because this sheet interleaves exhibitions/performances/local-life rows).
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")

$wsAll.Range("F2").Value  = 614
$wsAll.Range("F4").Value  = 1288
$wsAll.Range("F5").Value  = 1153
$wsAll.Range("F6").Value  = 14289
$wsAll.Range("F7").Value  = 16323
$wsAll.Range("F9").Value  = 90
$wsAll.Range("F18").Value = 99
$wsAll.Range("F20").Value = 1247
$wsAll.Range("F24").Value = 33
$wsAll.Range("F25").Value = 6572
$wsAll.Range("F27").Value = 15

# Row 31 ("【会员购严选】苏州·Come in joy动漫国潮文化节") is duplicated via
# Copy+Insert so formatting carries over, and its "想去人数" count is bumped to
# the updated value. This pushes every row below (归离之缘 / CF01 / 白日梦想 /
# 萤火 / 理想乡) down by one row, matching the row-count bump in <dimension>.
$wsAll.Rows(31).Copy()
$wsAll.Rows(31).Insert()
$wsAll.Range("F31").Value = 5707

# Column A's style didn't carry through the row-insert (Excel recomputed it
# from the row above), so restore it by pasting just the number/border format
# from the now-correctly-styled row below.
$wsAll.Range("A32").Copy()
$wsAll.Range("A31").PasteSpecial(-4122)

# The "萤火国潮文化节" row (now pushed down to row 36) also had its counter bumped.
$wsAll.Range("F36").Value = 4757
